$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.205.19"
$ws.Range("E2").Value = "  +3.29%  "
$ws.Range("D3").Value = "2.366.67"
$ws.Range("E3").Value = "  +1.45%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.91"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.60"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.609"
$ws.Range("E9").Value = "  -1.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.91"
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.45"
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("E13").Value = "  +1.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.976"
$ws.Range("E14").Value = "  -2.83%  "
$ws.Range("D15").Value = "2.727.43"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.22"
$ws.Range("E16").Value = "  -1.73%  "
$ws.Range("D17").Value = "2.374.75"
$ws.Range("E17").Value = "  +2.16%  "
$ws.Range("D18").Value = "45.187.66"
$ws.Range("E18").Value = "  +3.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.33"
$ws.Range("E19").Value = "  +9.60%  "
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.19"
$ws.Range("E21").Value = "  -4.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.04"
$ws.Range("E22").Value = "  -1.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.51"
$ws.Range("E23").Value = "  +1.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "258.37"
$ws.Range("E24").Value = "  -3.58%  "
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.06"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.21"
$ws.Range("E28").Value = "  -5.89%  "
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0972"
$ws.Range("E30").Value = "  +9.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.35"
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.14"
$ws.Range("E32").Value = "  -6.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "168.22"
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("E34").Value = "  +5.92%  "
$ws.Range("E35").Value = "  -1.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.116"
$ws.Range("E36").Value = "  +1.08%  "
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.96"
$ws.Range("E38").Value = "  +3.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0352"
$ws.Range("E39").Value = "  -3.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.89"
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.77"
$ws.Range("E41").Value = "  +3.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.12"
$ws.Range("E42").Value = "  -5.09%  "
$ws.Range("D43").Value = "1.895.36"
$ws.Range("E43").Value = "  +14.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "69.49"
$ws.Range("E44").Value = "  -3.05%  "
$ws.Range("E45").Value = "  -4.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.87"
$ws.Range("E46").Value = "  -5.96%  "
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "83.89"
$ws.Range("E48").Value = "  +9.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.64"
$ws.Range("E49").Value = "  +7.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.22"
$ws.Range("E50").Value = "  +2.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "110.05"
$ws.Range("E51").Value = "  -3.87%  "
